# Adds two new columns ("I0" in I, "IF" in J) to the sheet, mirroring the
# existing header style (copied from H1, which is bold/bordered/centered)
# and filling in the numeric data for rows 2-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers: copy H1's formatting onto I1/J1, then set the new labels ---
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-27 for columns I (I0) and J (IF) ---
$data = @(
    @(2, 5, 8),
    @(3, 1, 4),
    @(4, 1, 4),
    @(5, 2, 7),
    @(6, 1, 3),
    @(7, 1, 5),
    @(8, 1, 5),
    @(9, 1, 4),
    @(10, 1, 5),
    @(11, 1, 5),
    @(12, 1, 6),
    @(13, 1, 6),
    @(14, 1, 5),
    @(15, 1, 6),
    @(16, 1, 5),
    @(17, 1, 5),
    @(18, 1, 5),
    @(19, 1, 5),
    @(20, 1, 4),
    @(21, 1, 3),
    @(22, 6, 7),
    @(23, 7, 7),
    @(24, 8, 9),
    @(25, 8, 9),
    @(26, 3, 6),
    @(27, 9, 9)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
